$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modelIterations")

# Row 25 (Neural Networks / MLP - R^2)
$ws.Range("K25").Value = 62.9
$ws.Range("L25").Value = 46.2
$ws.Range("M25").Value = 39.200000000000003
$ws.Range("N25").Value = 22.7

# Row 26 (Neural Networks / MLP - RMSE)
$ws.Range("K26").Value = 84.4
$ws.Range("L26").Value = 64.400000000000006
$ws.Range("M26").Value = 58.3
$ws.Range("N26").Value = 32

# Row 27 (Neural Networks / MLP - MAE)
$ws.Range("K27").Value = 0.4
$ws.Range("L27").Value = 0.38
$ws.Range("M27").Value = 0.36
$ws.Range("N27").Value = 0.41

# Update the active selection to match the final cursor position
$ws.Range("K28").Select()
